$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell value updates per diff (cryptos list refresh)
$ws.Range("D2").Value = "42.412.63"
$ws.Range("E2").Value = "  +0.39%  "
$ws.Range("D3").Value = "2.283.93"
$ws.Range("E3").Value = "  -0.95%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "323.01"
$ws.Range("E5").Value = "  +2.07%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "102.34"
$ws.Range("E6").Value = "  -2.56%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.629"
$ws.Range("E7").Value = "  +0.17%  "
$ws.Range("E8").Value = "  +0.15%  "
$ws.Range("E9").Value = "  -0.70%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.63"
$ws.Range("E10").Value = "  -0.45%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0905"
$ws.Range("E11").Value = "  -0.33%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.30"
$ws.Range("E12").Value = "  -2.35%  "
$ws.Range("E13").Value = "  -0.74%  "
$ws.Range("E14").Value = "  -1.49%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.07"
$ws.Range("E15").Value = "  -2.55%  "
$ws.Range("D16").Value = "2.630.95"
$ws.Range("E16").Value = "  -0.92%  "
$ws.Range("D17").Value = "2.282.00"
$ws.Range("E17").Value = "  -0.65%  "
$ws.Range("D18").Value = "42.361.18"
$ws.Range("E18").Value = "  +0.51%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.35"
$ws.Range("E19").Value = "  -4.94%  "
$ws.Range("E20").Value = "  -0.64%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.70"
$ws.Range("E21").Value = "  +26.76%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.65"
$ws.Range("E22").Value = "  +2.99%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "73.00"
$ws.Range("E23").Value = "  -0.88%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "267.76"
$ws.Range("E24").Value = "  -4.53%  "
$ws.Range("E25").Value = "  -3.50%  "
$ws.Range("E26").Value = "  -0.13%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.83"
$ws.Range("E27").Value = "  -1.20%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.31"
$ws.Range("E28").Value = "  +2.03%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "22.43"
$ws.Range("E29").Value = "  -3.81%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "38.02"
$ws.Range("E30").Value = "  +5.80%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "164.31"
$ws.Range("E31").Value = "  -0.77%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.05"
$ws.Range("E32").Value = "  +2.61%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0875"
$ws.Range("E33").Value = "  -1.17%  "
$ws.Range("E34").Value = "  +0.83%  "
$ws.Range("E35").Value = "  -4.47%  "
$ws.Range("E36").Value = "  -12.99%  "
$ws.Range("E37").Value = "  -1.52%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0355"
$ws.Range("E38").Value = "  +0.40%  "
$ws.Range("E39").Value = "  +0.58%  "
$ws.Range("E40").Value = "  -6.86%  "
$ws.Range("E41").Value = "  +1.83%  "
$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  +0.02%  "
$ws.Range("B43").Value = "MultiversX"
$ws.Range("C43").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "68.41"
$ws.Range("E43").Value = "  -3.52%  "
$ws.Range("B44").Value = "Algorand"
$ws.Range("C44").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.225"
$ws.Range("E44").Value = "  -0.63%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "90.66"
$ws.Range("E45").Value = "  -10.84%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.22"
$ws.Range("E46").Value = "  +0.53%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "112.83"
$ws.Range("E47").Value = "  -3.49%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "79.91"
$ws.Range("E48").Value = "  +1.83%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.93"
$ws.Range("E49").Value = "  -1.86%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.23"
$ws.Range("E50").Value = "  -2.17%  "
$ws.Range("D51").Value = "1.591.01"
$ws.Range("E51").Value = "  +1.95%  "
